$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.726.98"
$ws.Range("E2").Value = "  +1.71%  "
$ws.Range("D3").Value = "2.091.14"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.386"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0845"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("E11").Value = "  +0.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.22"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.23%  "
$ws.Range("D13").Value = "2.400.39"
$ws.Range("E13").Value = "  -0.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.97"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.812"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.85%  "
$ws.Range("E16").Value = "  -0.55%  "
$ws.Range("D17").Value = "2.094.56"
$ws.Range("E17").Value = "  +1.75%  "
$ws.Range("D18").Value = "38.665.44"
$ws.Range("E18").Value = "  +2.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.74%  "
$ws.Range("D21").Value = "0.0₃0842"
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.45"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.45%  "
$ws.Range("E24").Value = "  -1.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.26%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.85"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.79%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.53"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.52%  "
$ws.Range("E28").Value = "  +5.94%  "
$ws.Range("E29").Value = "  +8.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.48"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.31%  "
$ws.Range("E32").Value = "  +0.36%  "
$ws.Range("E33").Value = "  +1.69%  "
$ws.Range("E34").Value = "  +0.66%  "
$ws.Range("E35").Value = "  +0.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.55"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.89%  "
$ws.Range("E37").Value = "  -0.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.58"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.01"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.56%  "
$ws.Range("E41").Value = "  +4.88%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.99"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.71%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.533.51"
$ws.Range("E43").Value = "  -0.80%  "
$ws.Range("E44").Value = "  -0.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0917"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.71"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.92%  "
$ws.Range("E47").Value = "  +1.05%  "
$ws.Range("E48").Value = "  -1.92%  "
$ws.Range("E49").Value = "  +0.93%  "
$ws.Range("E50").Value = "  -0.84%  "
$ws.Range("D51").Value = "2.286.74"
$ws.Range("E51").Value = "  -0.05%  "
